$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number-format on Price cells whose new values look numeric,
# so Excel stores them as literal text (matching the source data export)
# instead of silently converting to a Number type (which would also mangle
# values like "1.00" -> 1, "0.0000210" -> 2.1E-05, "3.40" -> 3.4, etc).
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D51').NumberFormat = "@"

$ws.Range('D2').Value2 = '96.996.46'
$ws.Range('E2').Value2 = '  +0.66%  '
$ws.Range('D3').Value2 = '3.687.24'
$ws.Range('E3').Value2 = '  +0.77%  '
$ws.Range('E4').Value2 = '  +0.03%  '
$ws.Range('D5').Value2 = '236.68'
$ws.Range('E5').Value2 = '  -2.14%  '
$ws.Range('D6').Value2 = '1.91'
$ws.Range('E6').Value2 = '  +0.66%  '
$ws.Range('D7').Value2 = '658.45'
$ws.Range('E7').Value2 = '  -0.03%  '
$ws.Range('E8').Value2 = '  -0.20%  '
$ws.Range('B9').Value2 = 'USDC'
$ws.Range('C9').Value2 = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D9').Value2 = '0.999'
$ws.Range('E9').Value2 = '  +0.00%  '
$ws.Range('B10').Value2 = 'Cardano'
$ws.Range('C10').Value2 = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D10').Value2 = '1.07'
$ws.Range('E10').Value2 = '  -2.14%  '
$ws.Range('D11').Value2 = '3.681.98'
$ws.Range('E11').Value2 = '  +0.67%  '
$ws.Range('D12').Value2 = '44.15'
$ws.Range('E12').Value2 = '  -1.60%  '
$ws.Range('E13').Value2 = '  +2.20%  '
$ws.Range('E14').Value2 = '  +10.64%  '
$ws.Range('D15').Value2 = '6.79'
$ws.Range('E15').Value2 = '  +1.75%  '
$ws.Range('B16').Value2 = 'WrappedBTC'
$ws.Range('C16').Value2 = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value2 = '96.731.67'
$ws.Range('E16').Value2 = '  +0.63%  '
$ws.Range('B17').Value2 = 'Polkadot'
$ws.Range('C17').Value2 = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D17').Value2 = '9.16'
$ws.Range('E17').Value2 = '  +3.10%  '
$ws.Range('B18').Value2 = 'WrappedEther'
$ws.Range('C18').Value2 = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value2 = '3.684.50'
$ws.Range('E18').Value2 = '  +0.81%  '
$ws.Range('B19').Value2 = 'Uniswap'
$ws.Range('C19').Value2 = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').Value2 = '13.03'
$ws.Range('E19').Value2 = '  +2.35%  '
$ws.Range('B20').Value2 = 'Chainlink'
$ws.Range('C20').Value2 = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D20').Value2 = '18.75'
$ws.Range('E20').Value2 = '  +2.74%  '
$ws.Range('B21').Value2 = 'Stellar'
$ws.Range('C21').Value2 = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D21').Value2 = '0.506'
$ws.Range('E21').Value2 = '  -5.19%  '
$ws.Range('B22').Value2 = 'BitcoinCash'
$ws.Range('C22').Value2 = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D22').Value2 = '518.34'
$ws.Range('E22').Value2 = '  -0.60%  '
$ws.Range('B23').Value2 = 'SuiNetwork'
$ws.Range('C23').Value2 = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D23').Value2 = '3.40'
$ws.Range('E23').Value2 = '  -1.26%  '
$ws.Range('B24').Value2 = 'PEPE'
$ws.Range('C24').Value2 = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D24').Value2 = '0.0000210'
$ws.Range('E24').Value2 = '  +2.89%  '
$ws.Range('B25').Value2 = 'NEARProtocol'
$ws.Range('C25').Value2 = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D25').Value2 = '6.94'
$ws.Range('E25').Value2 = '  +0.64%  '
$ws.Range('B26').Value2 = 'Hedera'
$ws.Range('C26').Value2 = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D26').Value2 = '0.204'
$ws.Range('E26').Value2 = '  +23.43%  '
$ws.Range('B27').Value2 = 'Litecoin'
$ws.Range('C27').Value2 = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D27').Value2 = '101.27'
$ws.Range('E27').Value2 = '  -1.02%  '
$ws.Range('B28').Value2 = 'Aptos'
$ws.Range('C28').Value2 = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D28').Value2 = '13.33'
$ws.Range('E28').Value2 = '  +2.77%  '
$ws.Range('B29').Value2 = 'InternetComputer(DFINITY)'
$ws.Range('C29').Value2 = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D29').Value2 = '12.48'
$ws.Range('E29').Value2 = '  +1.58%  '
$ws.Range('B30').Value2 = 'PancakeSwap'
$ws.Range('C30').Value2 = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').Value2 = '3.01'
$ws.Range('E30').Value2 = '  -0.15%  '
$ws.Range('B31').Value2 = 'Dai'
$ws.Range('C31').Value2 = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D31').Value2 = '0.999'
$ws.Range('E31').Value2 = '  -0.03%  '
$ws.Range('B32').Value2 = 'Cronos'
$ws.Range('C32').Value2 = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D32').Value2 = '0.191'
$ws.Range('E32').Value2 = '  +3.65%  '
$ws.Range('B33').Value2 = 'Fetch.AI'
$ws.Range('C33').Value2 = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D33').Value2 = '1.86'
$ws.Range('E33').Value2 = '  +1.94%  '
$ws.Range('B34').Value2 = 'Binance-PegBSC-USD'
$ws.Range('C34').Value2 = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D34').Value2 = '1.00'
$ws.Range('E34').Value2 = '  +0.19%  '
$ws.Range('B35').Value2 = 'EthereumClassic'
$ws.Range('C35').Value2 = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D35').Value2 = '32.13'
$ws.Range('E35').Value2 = '  -3.17%  '
$ws.Range('B36').Value2 = 'Bittensor'
$ws.Range('C36').Value2 = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D36').Value2 = '644.60'
$ws.Range('E36').Value2 = '  +3.07%  '
$ws.Range('B37').Value2 = 'PolygonEcosystemToken'
$ws.Range('C37').Value2 = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D37').Value2 = '0.592'
$ws.Range('E37').Value2 = '  +0.95%  '
$ws.Range('B38').Value2 = 'RenderToken'
$ws.Range('C38').Value2 = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D38').Value2 = '8.83'
$ws.Range('E38').Value2 = '  +1.23%  '
$ws.Range('B39').Value2 = 'USDe'
$ws.Range('C39').Value2 = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D39').Value2 = '1.00'
$ws.Range('E39').Value2 = '  +0.01%  '
$ws.Range('B40').Value2 = 'Algorand'
$ws.Range('C40').Value2 = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D40').Value2 = '0.504'
$ws.Range('E40').Value2 = '  +18.78%  '
$ws.Range('D41').Value2 = '6.87'
$ws.Range('E41').Value2 = '  +9.99%  '
$ws.Range('D42').Value2 = '2.07'
$ws.Range('E42').Value2 = '  +6.67%  '
$ws.Range('E43').Value2 = '  +1.26%  '
$ws.Range('D44').Value2 = '40.49'
$ws.Range('E44').Value2 = '  -12.87%  '
$ws.Range('B45').Value2 = 'ARBITRUM'
$ws.Range('C45').Value2 = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D45').Value2 = '0.959'
$ws.Range('E45').Value2 = '  +0.17%  '
$ws.Range('B46').Value2 = 'VeChain'
$ws.Range('C46').Value2 = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D46').Value2 = '0.0467'
$ws.Range('E46').Value2 = '  +3.28%  '
$ws.Range('B47').Value2 = 'WhiteBITCoin'
$ws.Range('C47').Value2 = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D47').Value2 = '23.63'
$ws.Range('E47').Value2 = '  +0.13%  '
$ws.Range('B48').Value2 = 'Stacks'
$ws.Range('C48').Value2 = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D48').Value2 = '2.28'
$ws.Range('E48').Value2 = '  -0.44%  '
$ws.Range('B49').Value2 = 'Cosmos'
$ws.Range('C49').Value2 = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D49').Value2 = '8.64'
$ws.Range('E49').Value2 = '  +1.47%  '
$ws.Range('B50').Value2 = 'MantraDAO'
$ws.Range('C50').Value2 = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range('D50').Value2 = '3.52'
$ws.Range('E50').Value2 = '  -1.63%  '
$ws.Range('B51').Value2 = 'OKB'
$ws.Range('C51').Value2 = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D51').Value2 = '53.78'
$ws.Range('E51').Value2 = '  -1.19%  '
